# revise and resubmit Democratization - Changes R&R Democratization
# Updates the "Model 1b" regression results table on Sheet1:
#   - Refreshes the coefficient stats for the existing two rows
#     ((Intercept), lrscale) with re-estimated values.
#   - Appends eight new predictor rows (age, educ, polint, sexMale,
#     surveyevs2008, surveywvs1994, surveywvs1999, surveywvs2005).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "icc" column (C) stores a numeric-looking value as literal TEXT in
# the source table. Force text entry (via a temporary "@" number format)
# so Excel doesn't auto-coerce the string into a Number, then drop the
# format back to Normal so no stray style index lingers on the cell.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$iccText = "0.0810811661067596"

# Row data: variable, coeff, icc(text), SE, n, name, lower, upper, n_country
$rows = @(
    @{ r = 2;  a = "(Intercept)";    b = 0.91;   d = 0.06774352886704478; e = 55488; g = 0.7985618950137113;  h = 1.021438104986289 },
    @{ r = 3;  a = "lrscale";        b = 0.389;  d = 0.01654970555267945; e = 55488; g = 0.3617757343658423;  h = 0.4162242656341577 },
    @{ r = 4;  a = "age";            b = -0.365; d = 0.02412527269832829; e = 55488; g = -0.40468607358875;   h = -0.32531392641125 },
    @{ r = 5;  a = "educ";           b = 0.292;  d = 0.01278016035775211; e = 55488; g = 0.2709766362114978;  h = 0.3130233637885022 },
    @{ r = 6;  a = "polint";         b = -0.299; d = 0.01479126976523935; e = 55488; g = -0.3233316387638187; h = -0.2746683612361813 },
    @{ r = 7;  a = "sexMale";        b = 0.014;  d = 0.008359596027554048; e = 55488; g = 0.0002484645346735904; h = 0.02775153546532641 },
    @{ r = 8;  a = "surveyevs2008";  b = -0.138; d = 0.01492129194477103; e = 55488; g = -0.1625455252491483; h = -0.1134544747508517 },
    @{ r = 9;  a = "surveywvs1994";  b = -0.094; d = 0.01500996242349729; e = 55488; g = -0.118691388186653;  h = -0.06930861181334697 },
    @{ r = 10; a = "surveywvs1999";  b = 0.193;  d = 0.0220727294712752;  e = 55488; g = 0.1566903600197523;  h = 0.2293096399802477 },
    @{ r = 11; a = "surveywvs2005";  b = 0.166;  d = 0.01924403826295021; e = 55488; g = 0.1343435570574469;  h = 0.1976564429425531 }
)

foreach ($row in $rows) {
    $r = $row.r

    $ws.Range("A$r").Value = $row.a
    $ws.Range("B$r").Value = $row.b
    Set-TextValue $ws.Range("C$r") $iccText
    $ws.Range("D$r").Value = $row.d
    $ws.Range("E$r").Value = $row.e
    $ws.Range("F$r").Value = "Model 1b"
    $ws.Range("G$r").Value = $row.g
    $ws.Range("H$r").Value = $row.h
    $ws.Range("I$r").Value = 20
}
